$wb = $excel.ActiveWorkbook

# --- Update ALUOp sheet view (zoom + selection) before adding the new sheet ---
$aluop = $wb.Worksheets.Item("ALUOp")
$aluop.Activate()
$excel.ActiveWindow.Zoom = 150
$aluop.Range("A9").Select()

# --- Add the new ALUFlag sheet after ALUOp ---
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "ALUFlag"

$ws.Columns.Item(1).ColumnWidth = 13.3
$ws.Columns.Item(2).ColumnWidth = 10.1

$full = $ws.Range("A1:B5")
$full.HorizontalAlignment = -4108
$full.Range("B1:B5").NumberFormat = "@"

$ws.Range("B1").Value = "ALUFlag"
$ws.Range("A1").Value = "Flag Symbol"
$ws.Range("A2").Value = "Z (Zero Flag)"
$ws.Range("A4").Value = "C (Carry Flag)"
$ws.Range("A3").Value = "S (Sign Flag)"
$ws.Range("B5").Value = "1000"
$ws.Range("A5").Value = "O (Overflow Flag)"

$ws.Range("B2").Value = "0001"
$ws.Range("B3").Value = "0010"
$ws.Range("B4").Value = "0100"

$ws.Range("A1:B1").Font.Bold = $true

# --- Activate ALUFlag and set its view state (becomes the active/selected tab) ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 210
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("A5").Select()

Write-Output "done"
